$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A422').Value = 'SE.COM.DURS'
$ws.Range('B422').Value = 'Compulsory education, duration (years)'
$ws.Range('F422').Value = 'Compulsory'
$ws.Range('G422').Value = 'Duration'
$ws.Range('A423').Value = 'SE.LPV.PRIM.MA'
$ws.Range('B423').Value = 'Learning poverty: Share of Male Children at the End-of-Primary age below minimum reading proficiency adjusted by Out-of-School Children (%)'
$ws.Range('F423').Value = 'Learning poverty'
$ws.Range('G423').Value = 'Primary education'
$ws.Range('A424').Value = 'SE.LPV.PRIM.FE'
$ws.Range('B424').Value = 'Learning poverty: Share of Female Children at the End-of-Primary age below minimum reading proficiency adjusted by Out-of-School Children (%)'
$ws.Range('F424').Value = 'Learning poverty'
$ws.Range('G424').Value = 'Primary education'
$ws.Range('A425').Value = 'SE.LPV.PRIM'
$ws.Range('B425').Value = 'Learning poverty: Share of Children at the End-of-Primary age below minimum reading proficiency adjusted by Out-of-School Children (%)'
$ws.Range('F425').Value = 'Learning poverty'
$ws.Range('G425').Value = 'Primary education'
$ws.Range('A426').Value = 'SE.LPV.PRIM.LD.FE'
$ws.Range('B426').Value = 'Female pupils below minimum reading proficiency at end of primary (%). Low GAML threshold'
$ws.Range('F426').Value = 'Learning poverty'
$ws.Range('G426').Value = 'Primary education'
$ws.Range('A427').Value = 'SE.LPV.PRIM.SD.FE'
$ws.Range('B427').Value = 'Female primary school age children out-of-school (%)'
$ws.Range('F427').Value = 'Learning poverty'
$ws.Range('G427').Value = 'Primary education'
$ws.Range('A428').Value = 'SE.LPV.PRIM.LD'
$ws.Range('B428').Value = 'Pupils below minimum reading proficiency at end of primary (%). Low GAML threshold'
$ws.Range('F428').Value = 'Learning poverty'
$ws.Range('G428').Value = 'Primary education'
$ws.Range('A429').Value = 'SE.LPV.PRIM.SD'
$ws.Range('B429').Value = 'Primary school age children out-of-school (%)'
$ws.Range('F429').Value = 'Learning poverty'
$ws.Range('G429').Value = 'Primary education'
$ws.Range('A430').Value = 'SE.LPV.PRIM.LD.MA'
$ws.Range('B430').Value = 'Male pupils below minimum reading proficiency at end of primary (%). Low GAML threshold'
$ws.Range('F430').Value = 'Learning poverty'
$ws.Range('A431').Value = 'SE.LPV.PRIM.SD.MA'
$ws.Range('B431').Value = 'Male primary school age children out-of-school (%)'
$ws.Range('F431').Value = 'Learning poverty'
$ws.Range('A432').Value = 'SE.ADT.LITR.ZS'
$ws.Range('B432').Value = 'Literacy rate, adult total (% of people ages 15 and above)'
$ws.Range('F432').Value = 'Literacy rate'
$ws.Range('G432').Value = 'Adult (ages 15 and above)'
$ws.Range('A433').Value = 'SE.ADT.LITR.MA.ZS'
$ws.Range('B433').Value = 'Literacy rate, adult male (% of males ages 15 and above)'
$ws.Range('F433').Value = 'Literacy rate'
$ws.Range('G433').Value = 'Adult (ages 15 and above)'
$ws.Range('A434').Value = 'SE.ADT.LITR.FE.ZS'
$ws.Range('B434').Value = 'Literacy rate, adult female (% of females ages 15 and above)'
$ws.Range('F434').Value = 'Literacy rate'
$ws.Range('G434').Value = 'Adult (ages 15 and above)'
$ws.Range('A435').Value = 'SE.ADT.1524.LT.FE.ZS'
$ws.Range('B435').Value = 'Literacy rate, youth female (% of females ages 15-24)'
$ws.Range('F435').Value = 'Literacy rate'
$ws.Range('G435').Value = 'Youth (ages 15-24)'
$ws.Range('A436').Value = 'SE.ADT.1524.LT.FM.ZS'
$ws.Range('B436').Value = 'Literacy rate, youth (ages 15-24), gender parity index (GPI)'
$ws.Range('F436').Value = 'Literacy rate'
$ws.Range('G436').Value = 'Youth (ages 15-24)'
$ws.Range('A437').Value = 'SE.ADT.1524.LT.MA.ZS'
$ws.Range('B437').Value = 'Literacy rate, youth male (% of males ages 15-24)'
$ws.Range('F437').Value = 'Literacy rate'
$ws.Range('G437').Value = 'Youth (ages 15-24)'
$ws.Range('A438').Value = 'SE.ADT.1524.LT.ZS'
$ws.Range('B438').Value = 'Literacy rate, youth total (% of people ages 15-24)'
$ws.Range('F438').Value = 'Literacy rate'
$ws.Range('G438').Value = 'Youth (ages 15-24)'
$ws.Range('A674').Value = 'FS.AST.CGOV.GD.ZS'
$ws.Range('B674').Value = 'Claims on central government, etc. (% GDP)'
$ws.Range('A675').Value = 'FS.AST.DOMS.GD.ZS'
$ws.Range('B675').Value = 'Domestic credit provided by financial sector (% of GDP)'
$ws.Range('G675').Value = 'Domestic'
$ws.Range('A676').Value = 'FS.AST.DOMO.GD.ZS'
$ws.Range('B676').Value = 'Claims on other sectors of the domestic economy (% of GDP)'
$ws.Range('G676').Value = 'Other domestic sectors'
$ws.Range('A677').Value = 'FS.AST.PRVT.GD.ZS'
$ws.Range('B677').Value = 'Domestic credit to private sector (% of GDP)'
$ws.Range('G677').Value = 'Private'
$ws.Range('A678').Value = 'FB.BNK.CAPA.ZS'
$ws.Range('B678').Value = 'Bank capital to assets ratio (%)'
$ws.Range('F678').Value = 'Bank'
$ws.Range('G678').Value = 'Capital to assets'
$ws.Range('F679').Value = 'Bank (miscellaneous)'
$ws.Range('A680').Value = 'FD.AST.PRVT.GD.ZS'
$ws.Range('B680').Value = 'Domestic credit to private sector by banks (% of GDP)'
$ws.Range('F680').Value = 'Deposit money banks'
$ws.Range('G680').Value = 'Private'
$ws.Range('A681').Value = 'FM.AST.CGOV.ZG.M3'
$ws.Range('B681').Value = 'Claims on central government (annual growth as % of broad money)'
$ws.Range('F681').Value = 'Monetary Survey'
$ws.Range('G681').Value = 'Central government'
$ws.Range('A682').Value = 'FM.AST.DOMS.CN'
$ws.Range('B682').Value = 'Net domestic credit (current LCU)'
$ws.Range('F682').Value = 'Monetary Survey'
$ws.Range('G682').Value = 'Domestic'
$ws.Range('A683').Value = 'FM.AST.NFRG.CN'
$ws.Range('B683').Value = 'Net foreign assets (current LCU)'
$ws.Range('F683').Value = 'Monetary Survey'
$ws.Range('G683').Value = 'Net foreign'
$ws.Range('A684').Value = 'FM.AST.DOMO.ZG.M3'
$ws.Range('B684').Value = 'Claims on other sectors of the domestic economy (annual growth as % of broad money)'
$ws.Range('F684').Value = 'Monetary Survey'
$ws.Range('G684').Value = 'Other domestic sectors'
$ws.Range('A685').Value = 'FM.AST.PRVT.GD.ZS'
$ws.Range('B685').Value = 'Monetary Sector credit to private sector (% GDP)'
$ws.Range('F685').Value = 'Monetary Survey'
$ws.Range('A686').Value = 'FM.AST.PRVT.ZG.M3'
$ws.Range('B686').Value = 'Claims on private sector (annual growth as % of broad money)'
$ws.Range('F686').Value = 'Monetary Survey'
$ws.Range('G686').Value = 'Private'
